$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 125.833336
$ws.Range("I9").Value = 131.75
$ws.Range("K9").Value = 131.75
$ws.Range("M9").Value = 37.25
$ws.Range("H17").Value = 1624.625
$ws.Range("J17").Value = 2000
$ws.Range("L17").Value = 6000
$ws.Range("N17").Value = -6336
$ws.Range("H39").Value = 99.90000000000001
$ws.Range("I39").Value = 99.90000000000001
$ws.Range("K39").Value = 299.7
$ws.Range("M39").Value = -3.700000000000045
$ws.Range("H100").Value = 1605.1538
$ws.Range("J100").Value = 3166
$ws.Range("L100").Value = 3166
$ws.Range("N100").Value = -4248
$ws.Range("H101").Value = 464
$ws.Range("I101").Value = 389.6
$ws.Range("K101").Value = 1168.8
$ws.Range("M101").Value = 453.1999999999998
$ws.Range("H106").Value = 22038.5
$ws.Range("I106").Value = 22038.5
$ws.Range("K106").Value = 22038.5
$ws.Range("M106").Value = -21407.5
$ws.Range("H125").Value = 989
$ws.Range("I125").Value = 989
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 8901
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -6441
$ws.Range("N125").Value = $null
$ws.Range("H137").Value = 3131.4666
$ws.Range("J137").Value = 5099.8
$ws.Range("L137").Value = 15299.4
$ws.Range("N137").Value = -20399.4
$ws.Range("H138").Value = 4478.175
$ws.Range("J138").Value = 4608.6055
$ws.Range("L138").Value = 13825.8165
$ws.Range("N138").Value = -24105.8165
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8170.7896
$ws.Range("I32").Value = 6190.8823
$ws.Range("K32").Value = 6190.8823
$ws.Range("M32").Value = -5903.8823
$ws.Range("H61").Value = 1362
$ws.Range("I61").Value = 1341.7142
$ws.Range("K61").Value = 1341.7142
$ws.Range("M61").Value = -1129.7142
$ws.Range("H97").Value = 1399.8
$ws.Range("J97").Value = 1244.5
$ws.Range("L97").Value = 1244.5
$ws.Range("N97").Value = -2236.5
$ws.Range("H102").Value = 806.6
$ws.Range("J102").Value = 2750
$ws.Range("L102").Value = 2750
$ws.Range("N102").Value = -5994
$ws.Range("H132").Value = 1544.2826
$ws.Range("I132").Value = 1456.6976
$ws.Range("K132").Value = 4370.0928
$ws.Range("M132").Value = -1840.0928
$ws.Range("H136").Value = 1362
$ws.Range("I136").Value = 1341.7142
$ws.Range("K136").Value = 4025.1426
$ws.Range("M136").Value = -1475.1426
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 599.6667
$ws.Range("I22").Value = 659.6
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 659.6
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -486.6
$ws.Range("N22").Value = -646
$ws.Range("H54").Value = 6000
$ws.Range("I54").Value = 6000
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 6000
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -5516
$ws.Range("N54").Value = $null
$ws.Range("H80").Value = 297.66666
$ws.Range("I80").Value = 426.5
$ws.Range("K80").Value = 426.5
$ws.Range("M80").Value = 571.5
$ws.Range("H83").Value = 297.66666
$ws.Range("I83").Value = 426.5
$ws.Range("K83").Value = 2132.5
$ws.Range("M83").Value = 2859.5
$ws.Range("H94").Value = 2904.8572
$ws.Range("I94").Value = 2387.3333
$ws.Range("K94").Value = 2387.3333
$ws.Range("M94").Value = -1936.3333
$ws.Range("H102").Value = 4727.5
$ws.Range("I102").Value = 4727.5
$ws.Range("K102").Value = 4727.5
$ws.Range("M102").Value = -1482.5
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = $null
$ws.Range("N134").Value = $null
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4922.2666
$ws.Range("I31").Value = 1994.7142
$ws.Range("J31").Value = 7483.875
$ws.Range("K31").Value = 1994.7142
$ws.Range("L31").Value = 7483.875
$ws.Range("M31").Value = -1699.7142
$ws.Range("N31").Value = -8073.875
$ws.Range("H34").Value = 4922.2666
$ws.Range("I34").Value = 1994.7142
$ws.Range("J34").Value = 7483.875
$ws.Range("K34").Value = 1994.7142
$ws.Range("L34").Value = 7483.875
$ws.Range("M34").Value = -1792.7142
$ws.Range("N34").Value = -7887.875
$ws.Range("H52").Value = 94600
$ws.Range("J52").Value = 94600
$ws.Range("L52").Value = 94600
$ws.Range("N52").Value = -95188
$ws.Range("H54").Value = 14000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 14000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 14000
$ws.Range("M54").Value = $null
$ws.Range("N54").Value = -15316
$ws.Range("H86").Value = 10565.667
$ws.Range("I86").Value = 9465
$ws.Range("K86").Value = 9465
$ws.Range("M86").Value = -8342
$ws.Range("H89").Value = 10565.667
$ws.Range("I89").Value = 9465
$ws.Range("K89").Value = 47325
$ws.Range("M89").Value = -41709
$ws.Range("H93").Value = 2290.4
$ws.Range("I93").Value = 2290.4
$ws.Range("K93").Value = 2290.4
$ws.Range("M93").Value = -418.4000000000001
$ws.Range("H105").Value = 1308.5
$ws.Range("I105").Value = 775.25
$ws.Range("K105").Value = 775.25
$ws.Range("M105").Value = 971.75
$ws.Range("H132").Value = 2614.158
$ws.Range("I132").Value = 1763.3334
$ws.Range("J132").Value = 4072.7144
$ws.Range("K132").Value = 5290.0002
$ws.Range("L132").Value = 12218.1432
$ws.Range("M132").Value = -2760.0002
$ws.Range("N132").Value = -17278.1432
$ws.Range("H141").Value = 19250
$ws.Range("J141").Value = 19250
$ws.Range("L141").Value = 19250
$ws.Range("N141").Value = -29610
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 110.55556
$ws.Range("I40").Value = 110.55556
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 442.22224
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -373.22224
$ws.Range("N40").Value = $null
$ws.Range("H42").Value = 1002
$ws.Range("J42").Value = 1002
$ws.Range("L42").Value = 3006
$ws.Range("N42").Value = -4074
$ws.Range("H131").Value = 1717.2858
$ws.Range("I131").Value = 1162.5
$ws.Range("J131").Value = 1939.2
$ws.Range("K131").Value = 3487.5
$ws.Range("L131").Value = 5817.6
$ws.Range("M131").Value = 1552.5
$ws.Range("N131").Value = -15897.6
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 7048
$ws.Range("I41").Value = 7048
$ws.Range("K41").Value = 7048
$ws.Range("M41").Value = -6693
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("M40").Value = -4864
$ws.Range("H61").Value = 5605.2144
$ws.Range("J61").Value = 4273.875
$ws.Range("L61").Value = 4273.875
$ws.Range("N61").Value = -4677.875
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = $null
$ws.Range("H100").Value = 3681.4167
$ws.Range("I100").Value = 1726.7142
$ws.Range("J100").Value = 6418
$ws.Range("K100").Value = 1726.7142
$ws.Range("L100").Value = 6418
$ws.Range("M100").Value = -1185.7142
$ws.Range("N100").Value = -7500
$ws.Range("H113").Value = 5605.2144
$ws.Range("J113").Value = 4273.875
$ws.Range("L113").Value = 4273.875
$ws.Range("N113").Value = -8613.875
$ws.Range("H122").Value = 9597.666999999999
$ws.Range("I122").Value = 9458.846
$ws.Range("K122").Value = 28376.538
$ws.Range("M122").Value = -25926.538
$ws.Range("H132").Value = 3853.394
$ws.Range("I132").Value = 3214.52
$ws.Range("K132").Value = 9643.559999999999
$ws.Range("M132").Value = -7113.559999999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 70076.5
$ws.Range("I51").Value = 70076
$ws.Range("K51").Value = 70076
$ws.Range("M51").Value = -69566
$ws.Range("H52").Value = 2050
$ws.Range("I52").Value = 2050
$ws.Range("K52").Value = 2050
$ws.Range("M52").Value = -1824
$ws.Range("H62").Value = 7258.684
$ws.Range("I62").Value = 5199.8
$ws.Range("J62").Value = 7994
$ws.Range("K62").Value = 5199.8
$ws.Range("L62").Value = 7994
$ws.Range("M62").Value = -4575.8
$ws.Range("N62").Value = -9242
$ws.Range("H65").Value = 7258.684
$ws.Range("I65").Value = 5199.8
$ws.Range("J65").Value = 7994
$ws.Range("K65").Value = 25999
$ws.Range("L65").Value = 39970
$ws.Range("M65").Value = -22879
$ws.Range("N65").Value = -46210
$ws.Range("H122").Value = 4408.857
$ws.Range("I122").Value = 4960.5835
$ws.Range("J122").Value = 1098.5
$ws.Range("K122").Value = 14881.7505
$ws.Range("L122").Value = 3295.5
$ws.Range("M122").Value = -12431.7505
$ws.Range("N122").Value = -8195.5
$ws.Range("H132").Value = 2020.4615
$ws.Range("I132").Value = 1809.5714
$ws.Range("J132").Value = 2266.5
$ws.Range("K132").Value = 5428.7142
$ws.Range("L132").Value = 6799.5
$ws.Range("M132").Value = -2898.7142
$ws.Range("N132").Value = -11859.5
$ws.Range("H136").Value = 2460.9312
$ws.Range("I136").Value = 781.1579
$ws.Range("K136").Value = 2343.4737
$ws.Range("M136").Value = 206.5263
